$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("PLAZAS_TARIFAS")
$ws2 = $wb.Worksheets.Item("PLAZAS_CANALES")
$ws3 = $wb.Worksheets.Item("TARIFAS_NACIONALES")

# --- Content corrections ---

# "CD JUAREZ " was an inconsistent label for the Ciudad Juarez plaza;
# standardize it (and the matching entry in PLAZAS_CANALES) to CIUDAD_JUAREZ.
$ws1.Range("A50").Value = "CIUDAD_JUAREZ"
$ws1.Range("A51").Value = "CIUDAD_JUAREZ"
$ws1.Range("A52").Value = "CIUDAD_JUAREZ"
$ws1.Range("A53").Value = "CIUDAD_JUAREZ"

# Row 67 was mislabeled QUERETARO; the 6/14/18/0-hour block it belongs to
# (rows 67-70) is actually SANLUIS, matching rows 68-70.
$ws1.Range("A67").Value = "SANLUIS"

# Keep PLAZAS_CANALES in sync with the corrected plaza name.
$ws2.Range("B6").Value = "CIUDAD_JUAREZ"

# --- Column width follow-up (content got wider -> best-fit column A) ---
$null = $ws1.Columns.Item(1).AutoFit()
$ws2.Columns.Item(7).ColumnWidth = 27.71
$ws2.Columns.Item(8).ColumnWidth = 20.46

# --- View-state / selection bookkeeping ---
$null = $ws1.Activate()
$excel.ActiveWindow.ScrollRow = 42
$excel.ActiveWindow.ScrollColumn = 1
$null = $ws1.Range("A54").Select()

$null = $ws3.Activate()
$null = $ws3.Range("F5").Select()

$null = $ws2.Activate()
$null = $ws2.Range("C7").Select()
